# Changement des polices par défaut pour les cellules vide (#29)
#
# 1) Switch the theme's font scheme ("super polices") so the major/minor
#    latin typeface is Calibri Light instead of Century Gothic.
# 2) Make the placeholder "Lorem" / "A vérifier" runs explicitly reference
#    the theme major-latin font (+mj-lt) instead of silently inheriting
#    whatever the old font scheme defined.
# 3) Shrink the autosized "A vérifier" textbox on slide 3 to fit the new
#    (narrower) font metrics.

$p = $ppt.ActivePresentation

# --- 1. Theme font scheme -------------------------------------------------
$fontScheme = $p.SlideMaster.Theme.ThemeFontScheme
$fontScheme.MajorFont.Latin = "Calibri Light"
$fontScheme.MinorFont.Latin = "Calibri Light"

# --- 2. Slide 1 : tag the "Lorem" placeholders with the major theme font --
$slide1 = $p.Slides.Item(1)

foreach ($name in @("nomBatiment", "adresse", "dateDeConstruction", "surfaceTotaleChauffe", "dateDeRenovation", "Elements de contexte sur le bâtiment")) {
    $slide1.Shapes.Item($name).TextFrame.TextRange.Font.Name = "+mj-lt"
}

# --- 3. Slide 3 : "A vérifier" textbox ------------------------------------
$slide3 = $p.Slides.Item(3)
$averifier = $slide3.Shapes.Item("ZoneTexte 2")
$averifier.TextFrame.TextRange.Font.Name = "+mj-lt"

# Resize the autosized textbox to the new (narrower) width computed by
# PowerPoint after switching the run to the major theme font.
$averifier.Width = 928075 / 914400 * 72
